# Rename header row cells: replace underscores with spaces for the
# specific columns changed in the commit "Version Final 03 fecha: 25/05/2023"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Fecha y hora presentacion"
$ws.Range("D1").Value = "Num de Operacion"
$ws.Range("E1").Value = "Periodo de declaracion"
$ws.Range("G1").Value = "total a pagar"
$ws.Range("H1").Value = "Vigente hasta"
$ws.Range("I1").Value = "Linea de Captura"
$ws.Range("K1").Value = "Impuesto a favor"
